# Sprint 2 backlog update
# Product backlog and User stories have been updated to reflect the work
# finished in Sprint 2: statuses moved forward, finish dates filled in,
# a couple of assignees trimmed, one user-story wording tweak, and one
# user story retitled/rescoped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "search based on location input" (touched first so that the
#     newly-introduced "7th Feb" value lands in the shared string table
#     ahead of the other brand-new strings, matching how the sheet was
#     actually edited) ---
$ws.Range("F6").Value = "5th Feb"
$ws.Range("G6").Value = "7th Feb"
$ws.Range("K6").Value = "Completed"

# --- Row 24: user story retitled ---
$ws.Range("C24").Value = "As a programmer, I want to host the website on the server"

# --- Row 12: price information on a marker (was "placeholder") ---
$ws.Range("C12").Value = "As a user, I can see the price information on a marker for each result on the map"

# --- Row 5: "search based on live location" ---
$ws.Range("F5").Value = "3rd Feb"
$ws.Range("K5").Value = "In Progress"

# --- Row 11: "view the search results on a map" ---
$ws.Range("G11").Value = "6th Feb"
$ws.Range("K11").Value = "Completed"

# --- Row 15: "set a price range for the search" ---
$ws.Range("G15").Value = "31st Jan"
$ws.Range("K15").Value = "Completed"

# --- Row 18: "navigate with a keyboard" ---
$ws.Range("F18").Value = "6th Feb"
$ws.Range("G18").Value = "6th Feb"
$ws.Range("K18").Value = "Completed"

# --- Row 19: "view the data regardless of colour blindness" ---
$ws.Range("F19").Value = "6th Feb"
$ws.Range("G19").Value = "6th Feb"
$ws.Range("K19").Value = "Completed"

# --- Row 21: "integrate the database, backend and frontend" ---
$ws.Range("G21").Value = "6th Feb"
$ws.Range("K21").Value = "Completed"

# --- Row 22: "store longtitude/latitude coordinates" ---
$ws.Range("E22").Value = "Kamila"
$ws.Range("F22").Value = "3rd Feb"
$ws.Range("G22").Value = "5th Feb"
$ws.Range("K22").Value = "Completed"

# --- Row 23: "optimise a webpage based on user feedback" ---
$ws.Range("G23").Value = "7th Feb"
$ws.Range("K23").Value = "Completed"

# --- Row 25: "create unit tests for all existing functionality" ---
$ws.Range("E25").Value = "Stas"
$ws.Range("F25").Value = "3rd Feb"
$ws.Range("G25").Value = "7th Feb"
$ws.Range("K25").Value = "Completed"

# Update the active selection to reflect where the author left off working
$ws.Range("H18").Select()
